# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: a handful of Price values are plain decimal numbers (e.g. "1.001",
# "241.80") that Excel would otherwise auto-convert to the number type,
# silently dropping the original text formatting (trailing zeros etc.).
# Force those particular cells to Text format first so the literal string
# is preserved, matching the source data exactly.
$textCells = @("D4","D5","D6","D8","D11","D13","D15","D16","D17","D19","D21","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D40","D41","D42","D43","D44","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.304.61"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.873.07"

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.001"

# Row 5 - XRP
$ws.Range("D5").Value = "0.7079"
$ws.Range("E5").Value = "  -0.56%  "

# Row 6 - BNB
$ws.Range("D6").Value = "241.80"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.07785"
$ws.Range("E8").Value = "  +0.98%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.12%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -1.15%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.08389"
$ws.Range("E11").Value = "  +0.14%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.880.73"
$ws.Range("E12").Value = "  +0.05%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "5.235"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +0.19%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "91.03"
$ws.Range("E15").Value = "  -0.74%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "6.137"
$ws.Range("E16").Value = "  +2.70%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.000008344"
$ws.Range("E17").Value = "  +0.91%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "29.316.00"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "239.96"
$ws.Range("E19").Value = "  -1.70%  "

# Row 20 - now WrappedliquidstakedEther2.0 (was Avalanche)
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.130.59"
$ws.Range("E20").Value = "  -0.39%  "

# Row 21 - now Avalanche (was WrappedliquidstakedEther2.0)
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "13.20"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.08%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "7.749"
$ws.Range("E23").Value = "  -1.91%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  +0.08%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "0.1591"
$ws.Range("E25").Value = "  -1.72%  "

# Row 26 - Monero
$ws.Range("D26").Value = "162.70"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "9.027"
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  -0.57%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  -0.26%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "4.414"
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "4.329"
$ws.Range("E31").Value = "  +0.11%  "

# Row 32 - Toncoin
$ws.Range("D32").Value = "1.242"
$ws.Range("E32").Value = "  -3.89%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.05343"
$ws.Range("E33").Value = "  +2.26%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "0.7511"
$ws.Range("E35").Value = "  -3.31%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  -0.09%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.243.60"
$ws.Range("E39").Value = "  +6.72%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "2.731"
$ws.Range("E40").Value = "  +0.35%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "6.532"
$ws.Range("E41").Value = "  +1.80%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.8925"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43 - Quant
$ws.Range("D43").Value = "109.34"
$ws.Range("E43").Value = "  +4.57%  "

# Row 44 - Aave
$ws.Range("D44").Value = "72.32"
$ws.Range("E44").Value = "  -1.55%  "

# Row 45 - PaxDollar
$ws.Range("E45").Value = "  +0.05%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +5.03%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "2.017.59"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "0.5201"
$ws.Range("E48").Value = "  -0.05%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  -0.46%  "

# Row 50 - EnergySwap
$ws.Range("E50").Value = "  +0.54%  "

# Row 51 - TheSandbox
$ws.Range("D51").Value = "0.4342"
$ws.Range("E51").Value = "  +0.73%  "
